# unitdata_TYNDP-2020.xlsx edit:
#  - Change the AutoFilter criteria on column B (Generator_ID) from
#    {Other non-RES, Other non-RES P, Other RES} to {Nuclear}. The Country
#    filter already active on column A (FI00) is left untouched; Excel
#    recomputes each row's Hidden state from the combination of both
#    filters.
#  - Select column I (whole column) as the new active selection, matching
#    the saved cursor/selection position in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:H4457")
[void]$dataRange.AutoFilter(2, @("Nuclear"), 7)

[void]$ws.Columns("I").Select()
